# Refresh cryptos list figures (price / 1h volume change) per latest scrape.
# Numeric-looking price strings are written with a leading quote so Excel
# keeps them as text (matching the original inline-string cell type) instead
# of silently converting them to numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.997.59'
$ws.Cells.Item(2, 5).Value = '  -0.43%  '
$ws.Cells.Item(3, 4).Value = '2.790.08'
$ws.Cells.Item(3, 5).Value = '  -1.78%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).Value = "'359.43"
$ws.Cells.Item(5, 5).Value = '  -0.41%  '
$ws.Cells.Item(6, 4).Value = "'109.71"
$ws.Cells.Item(6, 5).Value = '  -3.20%  '
$ws.Cells.Item(7, 5).Value = '  -3.04%  '
$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 5).Value = '  -2.68%  '
$ws.Cells.Item(10, 4).Value = "'40.21"
$ws.Cells.Item(10, 5).Value = '  -3.35%  '
$ws.Cells.Item(11, 4).Value = "'0.0850"
$ws.Cells.Item(11, 5).Value = '  -1.73%  '
$ws.Cells.Item(12, 4).Value = "'0.134"
$ws.Cells.Item(12, 5).Value = '  +1.48%  '
$ws.Cells.Item(13, 4).Value = "'19.50"
$ws.Cells.Item(13, 5).Value = '  -2.71%  '
$ws.Cells.Item(14, 5).Value = '  -3.16%  '
$ws.Cells.Item(15, 4).Value = '3.232.61'
$ws.Cells.Item(15, 5).Value = '  -1.74%  '
$ws.Cells.Item(16, 4).Value = '2.801.18'
$ws.Cells.Item(16, 5).Value = '  -0.41%  '
$ws.Cells.Item(17, 4).Value = "'0.943"
$ws.Cells.Item(17, 5).Value = '  +3.64%  '
$ws.Cells.Item(18, 4).Value = '51.943.83'
$ws.Cells.Item(18, 5).Value = '  -0.32%  '
$ws.Cells.Item(19, 5).Value = '  -1.91%  '
$ws.Cells.Item(20, 4).Value = "'3.11"
$ws.Cells.Item(20, 5).Value = '  -1.46%  '
$ws.Cells.Item(21, 4).Value = "'13.15"
$ws.Cells.Item(21, 5).Value = '  -2.95%  '
$ws.Cells.Item(22, 5).Value = '  -2.00%  '
$ws.Cells.Item(23, 4).Value = "'70.31"
$ws.Cells.Item(23, 5).Value = '  -0.18%  '
$ws.Cells.Item(24, 4).Value = "'270.64"
$ws.Cells.Item(24, 5).Value = '  +0.85%  '
$ws.Cells.Item(25, 4).Value = "'2.76"
$ws.Cells.Item(25, 5).Value = '  -2.71%  '
$ws.Cells.Item(26, 4).Value = "'26.52"
$ws.Cells.Item(27, 5).Value = '  -0.03%  '
$ws.Cells.Item(28, 4).Value = "'0.162"
$ws.Cells.Item(28, 5).Value = '  +15.31%  '
$ws.Cells.Item(29, 4).Value = "'10.32"
$ws.Cells.Item(29, 5).Value = '  -1.10%  '
$ws.Cells.Item(30, 4).Value = "'2.29"
$ws.Cells.Item(30, 5).Value = '  +1.53%  '
$ws.Cells.Item(31, 4).Value = "'0.0474"
$ws.Cells.Item(31, 5).Value = '  -2.33%  '
$ws.Cells.Item(32, 4).Value = "'52.01"
$ws.Cells.Item(32, 5).Value = '  -3.68%  '
$ws.Cells.Item(33, 4).Value = "'34.33"
$ws.Cells.Item(33, 5).Value = '  -1.21%  '
$ws.Cells.Item(34, 4).Value = "'5.75"
$ws.Cells.Item(34, 5).Value = '  -2.41%  '
$ws.Cells.Item(35, 5).Value = '  +0.24%  '
$ws.Cells.Item(36, 5).Value = '  -4.12%  '
$ws.Cells.Item(37, 5).Value = '  -0.06%  '
$ws.Cells.Item(38, 4).Value = "'19.08"
$ws.Cells.Item(38, 5).Value = '  +3.59%  '
$ws.Cells.Item(39, 4).Value = "'3.20"
$ws.Cells.Item(39, 5).Value = '  -2.22%  '
$ws.Cells.Item(40, 5).Value = '  -3.92%  '
$ws.Cells.Item(41, 4).Value = "'2.65"
$ws.Cells.Item(41, 5).Value = '  +3.48%  '
$ws.Cells.Item(42, 5).Value = '  -2.18%  '
$ws.Cells.Item(43, 5).Value = '  -1.03%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = "'21.99"
$ws.Cells.Item(44, 5).Value = '  -8.55%  '
$ws.Cells.Item(45, 2).Value = 'Monero'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(45, 4).Value = "'119.50"
$ws.Cells.Item(45, 5).Value = '  -6.64%  '
$ws.Cells.Item(46, 4).Value = '2.081.61'
$ws.Cells.Item(46, 5).Value = '  -1.52%  '
$ws.Cells.Item(47, 5).Value = '  -4.79%  '
$ws.Cells.Item(49, 5).Value = '  -0.83%  '
$ws.Cells.Item(50, 4).Value = "'0.955"
$ws.Cells.Item(50, 5).Value = '  -5.07%  '
$ws.Cells.Item(51, 5).Value = '  -2.97%  '
